# Apply the "augRie radial Dam break" GPU/CUDA measurements alongside the
# existing CPU measurements (columns E/F mirror columns A/B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header for the second measurement block (col E/F) ---
$ws.Range("E1").Value = "Problem; 2512x2512 augRie, radial Dam break"

# --- Single-node CPU / CUDA ---
$ws.Range("E4").Value = "Single-node CPU"
$ws.Range("F4").Value = 345

$ws.Range("E5").Value = "Single-node CUDA"
$ws.Range("F5").Value = 46

# --- MPI CPU / CUDA ---
$ws.Range("E6").Value = "MPI CPU"
$ws.Range("F6").Value = 465

$ws.Range("E7").Value = "MPI CUDA"
$ws.Range("F7").Value = 29

# --- STARPU_LOCAL_WORK_STEALING 2x2 (mirrors A10/B10's row) ---
$ws.Range("E10").Value = "STARPU_LOCAL_WORK_STEALING 2x2"
$ws.Range("F10").Value = 677

# --- dm 2x2 / dmda 2x2 ---
$ws.Range("E15").Value = "dm 2x2"
$ws.Range("F15").Value = 26

$ws.Range("E16").Value = "dmda 2x2"
$ws.Range("F16").Value = 26

# --- Column E got wider to fit the longer labels ---
$ws.Columns.Item(5).ColumnWidth = 37.46

# --- View: scrolled down/right a bit, new active cell ---
$ws.Range("G12").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 2
